$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(38, 8).Value = 535.6667
$ws.Cells.Item(38, 9).Value = 535.6667
$ws.Cells.Item(38, 10).Value = 0
$ws.Cells.Item(38, 11).Value = 1607.0001
$ws.Cells.Item(38, 12).Value = 0
$ws.Cells.Item(38, 13).Value = -1235.0001
$ws.Cells.Item(38, 14).ClearContents()
$ws.Cells.Item(64, 8).Value = 2819.5107
$ws.Cells.Item(64, 9).Value = 2613.8572
$ws.Cells.Item(64, 10).Value = 3122.5789
$ws.Cells.Item(64, 11).Value = 2613.8572
$ws.Cells.Item(64, 12).Value = 3122.5789
$ws.Cells.Item(64, 13).Value = -2365.8572
$ws.Cells.Item(64, 14).Value = -3618.5789
$ws.Cells.Item(67, 8).Value = 2819.5107
$ws.Cells.Item(67, 9).Value = 2613.8572
$ws.Cells.Item(67, 10).Value = 3122.5789
$ws.Cells.Item(67, 11).Value = 2613.8572
$ws.Cells.Item(67, 12).Value = 3122.5789
$ws.Cells.Item(67, 13).Value = -1755.8572
$ws.Cells.Item(67, 14).Value = -4838.5789
$ws.Cells.Item(86, 8).Value = 5286.6
$ws.Cells.Item(86, 9).Value = 4599.1
$ws.Cells.Item(86, 10).Value = 6661.6
$ws.Cells.Item(86, 11).Value = 4599.1
$ws.Cells.Item(86, 12).Value = 6661.6
$ws.Cells.Item(86, 13).Value = -3476.1
$ws.Cells.Item(86, 14).Value = -8907.6
$ws.Cells.Item(89, 8).Value = 5286.6
$ws.Cells.Item(89, 9).Value = 4599.1
$ws.Cells.Item(89, 10).Value = 6661.6
$ws.Cells.Item(89, 11).Value = 22995.5
$ws.Cells.Item(89, 12).Value = 33308
$ws.Cells.Item(89, 13).Value = -17379.5
$ws.Cells.Item(89, 14).Value = -44540
$ws.Cells.Item(129, 8).Value = 1016.4286
$ws.Cells.Item(129, 9).Value = 0
$ws.Cells.Item(129, 10).Value = 1016.4286
$ws.Cells.Item(129, 11).Value = 0
$ws.Cells.Item(129, 12).Value = 3049.2858
$ws.Cells.Item(129, 13).ClearContents()
$ws.Cells.Item(129, 14).Value = -13049.2858
$ws.Cells.Item(138, 8).Value = 2742146.5
$ws.Cells.Item(138, 9).Value = 1040.6666
$ws.Cells.Item(138, 10).Value = 5409168.5
$ws.Cells.Item(138, 11).Value = 3121.9998
$ws.Cells.Item(138, 12).Value = 16227505.5
$ws.Cells.Item(138, 13).Value = 2018.0002
$ws.Cells.Item(138, 14).Value = -16237785.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5146.9185
$ws.Cells.Item(32, 9).Value = 3023.9307
$ws.Cells.Item(32, 10).Value = 16065.143
$ws.Cells.Item(32, 11).Value = 3023.9307
$ws.Cells.Item(32, 12).Value = 16065.143
$ws.Cells.Item(32, 13).Value = -2736.9307
$ws.Cells.Item(32, 14).Value = -16639.143
$ws.Cells.Item(63, 8).Value = 3117.9412
$ws.Cells.Item(63, 9).Value = 2700.8333
$ws.Cells.Item(63, 10).Value = 3345.4546
$ws.Cells.Item(63, 11).Value = 2700.8333
$ws.Cells.Item(63, 12).Value = 3345.4546
$ws.Cells.Item(63, 13).Value = -2014.8333
$ws.Cells.Item(63, 14).Value = -4717.4546
$ws.Cells.Item(66, 8).Value = 3117.9412
$ws.Cells.Item(66, 9).Value = 2700.8333
$ws.Cells.Item(66, 10).Value = 3345.4546
$ws.Cells.Item(66, 11).Value = 13504.1665
$ws.Cells.Item(66, 12).Value = 16727.273
$ws.Cells.Item(66, 13).Value = -10072.1665
$ws.Cells.Item(66, 14).Value = -23591.273
$ws.Cells.Item(74, 8).Value = 27145.104
$ws.Cells.Item(74, 9).Value = 32354.906
$ws.Cells.Item(74, 10).Value = 3328.8572
$ws.Cells.Item(74, 11).Value = 32354.906
$ws.Cells.Item(74, 12).Value = 3328.8572
$ws.Cells.Item(74, 13).Value = -31480.906
$ws.Cells.Item(74, 14).Value = -5076.8572
$ws.Cells.Item(75, 8).Value = 0
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 12).Value = 0
$ws.Cells.Item(75, 14).ClearContents()
$ws.Cells.Item(77, 8).Value = 27145.104
$ws.Cells.Item(77, 9).Value = 32354.906
$ws.Cells.Item(77, 10).Value = 3328.8572
$ws.Cells.Item(77, 11).Value = 161774.53
$ws.Cells.Item(77, 12).Value = 16644.286
$ws.Cells.Item(77, 13).Value = -157406.53
$ws.Cells.Item(77, 14).Value = -25380.286
$ws.Cells.Item(78, 8).Value = 0
$ws.Cells.Item(78, 10).Value = 0
$ws.Cells.Item(78, 12).Value = 0
$ws.Cells.Item(78, 14).ClearContents()
$ws.Cells.Item(119, 8).Value = 30090.908
$ws.Cells.Item(119, 10).Value = 30090.908
$ws.Cells.Item(119, 12).Value = 30090.908
$ws.Cells.Item(119, 14).Value = -39766.908
$ws.Cells.Item(132, 8).Value = 2091.65
$ws.Cells.Item(132, 9).Value = 1559.5454
$ws.Cells.Item(132, 10).Value = 2742
$ws.Cells.Item(132, 11).Value = 4678.6362
$ws.Cells.Item(132, 12).Value = 8226
$ws.Cells.Item(132, 13).Value = -2148.6362
$ws.Cells.Item(132, 14).Value = -13286
$ws.Cells.Item(137, 8).Value = 64250
$ws.Cells.Item(137, 10).Value = 64250
$ws.Cells.Item(137, 12).Value = 64250
$ws.Cells.Item(137, 14).Value = -74450

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(51, 8).Value = 29649.5
$ws.Cells.Item(51, 10).Value = 29649.5
$ws.Cells.Item(51, 12).Value = 29649.5
$ws.Cells.Item(51, 14).Value = -30631.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 8930227
$ws.Cells.Item(31, 9).Value = 1126.4595
$ws.Cells.Item(31, 10).Value = 26318474
$ws.Cells.Item(31, 11).Value = 1126.4595
$ws.Cells.Item(31, 12).Value = 26318474
$ws.Cells.Item(31, 13).Value = -831.4594999999999
$ws.Cells.Item(31, 14).Value = -26319064
$ws.Cells.Item(34, 8).Value = 8930227
$ws.Cells.Item(34, 9).Value = 1126.4595
$ws.Cells.Item(34, 10).Value = 26318474
$ws.Cells.Item(34, 11).Value = 1126.4595
$ws.Cells.Item(34, 12).Value = 26318474
$ws.Cells.Item(34, 13).Value = -924.4594999999999
$ws.Cells.Item(34, 14).Value = -26318878
$ws.Cells.Item(58, 8).Value = 3423.025
$ws.Cells.Item(58, 9).Value = 3697.3428
$ws.Cells.Item(58, 10).Value = 1502.8
$ws.Cells.Item(58, 11).Value = 3697.3428
$ws.Cells.Item(58, 12).Value = 1502.8
$ws.Cells.Item(58, 13).Value = -3494.3428
$ws.Cells.Item(58, 14).Value = -1908.8
$ws.Cells.Item(100, 8).Value = 31548
$ws.Cells.Item(100, 10).Value = 31548
$ws.Cells.Item(100, 12).Value = 31548
$ws.Cells.Item(100, 14).Value = -33712
$ws.Cells.Item(134, 8).Value = 1742.6865
$ws.Cells.Item(134, 9).Value = 1697.9108
$ws.Cells.Item(134, 10).Value = 1970.6364
$ws.Cells.Item(134, 11).Value = 5093.732400000001
$ws.Cells.Item(134, 12).Value = 5911.9092
$ws.Cells.Item(134, 13).Value = -2558.732400000001
$ws.Cells.Item(134, 14).Value = -10981.9092
$ws.Cells.Item(136, 8).Value = 3423.025
$ws.Cells.Item(136, 9).Value = 3697.3428
$ws.Cells.Item(136, 10).Value = 1502.8
$ws.Cells.Item(136, 11).Value = 11092.0284
$ws.Cells.Item(136, 12).Value = 4508.4
$ws.Cells.Item(136, 13).Value = -8542.028399999999
$ws.Cells.Item(136, 14).Value = -9608.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(49, 8).Value = 817.6667
$ws.Cells.Item(109, 8).Value = 2173.4
$ws.Cells.Item(109, 9).Value = 1107.1
$ws.Cells.Item(109, 10).Value = 3239.7
$ws.Cells.Item(109, 11).Value = 3321.3
$ws.Cells.Item(109, 12).Value = 9719.099999999999
$ws.Cells.Item(109, 13).Value = -2281.3
$ws.Cells.Item(109, 14).Value = -11799.1
$ws.Cells.Item(115, 8).Value = 2142.4
$ws.Cells.Item(115, 9).Value = 1021
$ws.Cells.Item(115, 10).Value = 2890
$ws.Cells.Item(115, 11).Value = 3063
$ws.Cells.Item(115, 12).Value = 8670
$ws.Cells.Item(115, 13).Value = -1888
$ws.Cells.Item(115, 14).Value = -11020
$ws.Cells.Item(131, 8).Value = 918.1900000000001
$ws.Cells.Item(131, 10).Value = 920.2959
$ws.Cells.Item(131, 12).Value = 2760.8877
$ws.Cells.Item(131, 14).Value = -12840.8877

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 115084.336
$ws.Cells.Item(122, 9).Value = 145473.42
$ws.Cells.Item(122, 10).Value = 8722.5
$ws.Cells.Item(122, 11).Value = 436420.26
$ws.Cells.Item(122, 12).Value = 26167.5
$ws.Cells.Item(122, 13).Value = -433970.26
$ws.Cells.Item(122, 14).Value = -31067.5
$ws.Cells.Item(132, 8).Value = 3032969.5
$ws.Cells.Item(132, 9).Value = 2533.8518
$ws.Cells.Item(132, 10).Value = 16669930
$ws.Cells.Item(132, 11).Value = 7601.555399999999
$ws.Cells.Item(132, 12).Value = 50009790
$ws.Cells.Item(132, 13).Value = -5071.555399999999
$ws.Cells.Item(132, 14).Value = -50014850

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 861.9167
$ws.Cells.Item(46, 9).Value = 909.6
$ws.Cells.Item(46, 10).Value = 827.8570999999999
$ws.Cells.Item(46, 11).Value = 909.6
$ws.Cells.Item(46, 12).Value = 827.8570999999999
$ws.Cells.Item(46, 13).Value = -721.6
$ws.Cells.Item(46, 14).Value = -1203.8571
$ws.Cells.Item(64, 8).Value = 33500
$ws.Cells.Item(64, 10).Value = 33500
$ws.Cells.Item(64, 12).Value = 33500
$ws.Cells.Item(64, 14).Value = -33950
$ws.Cells.Item(67, 8).Value = 33500
$ws.Cells.Item(67, 10).Value = 33500
$ws.Cells.Item(67, 12).Value = 33500
$ws.Cells.Item(67, 14).Value = -35060

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(57, 8).Value = 24631.578
$ws.Cells.Item(57, 10).Value = 24631.578
$ws.Cells.Item(57, 12).Value = 24631.578
$ws.Cells.Item(57, 14).Value = -26139.578
$ws.Cells.Item(64, 8).Value = 29500
$ws.Cells.Item(64, 10).Value = 29500
$ws.Cells.Item(64, 12).Value = 29500
$ws.Cells.Item(64, 14).Value = -29996
$ws.Cells.Item(67, 8).Value = 29500
$ws.Cells.Item(67, 10).Value = 29500
$ws.Cells.Item(67, 12).Value = 29500
$ws.Cells.Item(67, 14).Value = -31216
$ws.Cells.Item(119, 8).Value = 30729.75
$ws.Cells.Item(119, 10).Value = 30729.75
$ws.Cells.Item(119, 12).Value = 30729.75
$ws.Cells.Item(119, 14).Value = -40405.75
$ws.Cells.Item(126, 8).Value = 66667668
$ws.Cells.Item(126, 9).Value = 854
$ws.Cells.Item(126, 10).Value = 200001300
$ws.Cells.Item(126, 11).Value = 2562
$ws.Cells.Item(126, 12).Value = 600003900
$ws.Cells.Item(126, 13).Value = -92
$ws.Cells.Item(126, 14).Value = -600008840
$ws.Cells.Item(132, 8).Value = 2481.1875
$ws.Cells.Item(132, 9).Value = 2518.4167
$ws.Cells.Item(132, 10).Value = 2369.5
$ws.Cells.Item(132, 11).Value = 7555.250100000001
$ws.Cells.Item(132, 12).Value = 7108.5
$ws.Cells.Item(132, 13).Value = -5025.250100000001
$ws.Cells.Item(132, 14).Value = -12168.5
$ws.Cells.Item(136, 8).Value = 1602.869
$ws.Cells.Item(136, 9).Value = 1408.8474
$ws.Cells.Item(136, 10).Value = 2060.76
$ws.Cells.Item(136, 11).Value = 4226.5422
$ws.Cells.Item(136, 12).Value = 6182.280000000001
$ws.Cells.Item(136, 13).Value = -1676.5422
$ws.Cells.Item(136, 14).Value = -11282.28
